# Commit: "Modified Kiel presentation: #3"
#
# The canonical-XML diff for this commit only touches the slide-transition
# <mc:AlternateContent>/<mc:Choice>/<mc:Fallback> wrapper on 4 slides
# (the "morph" transition on slide 10, and the "p14:dur" timed transitions
# on slides 11-13). In every hunk the namespace declaration for the
# version-gated prefix (p159:/p14:) is simply hoisted from <mc:Choice> up
# onto <mc:AlternateContent> (with a matching xmlns="" reset added on the
# sibling <mc:Fallback>) - every transition attribute/value/child already
# present (spd, advClick, advTm, p14:dur, the morph/fade effect children...)
# is left byte-for-byte identical. It is a namespace-serialization
# normalization with no visual or behavioral effect: the slide transitions
# keep exactly the same speed/advance/effect settings before and after.
#
# Re-assert the existing transition settings on the affected slides so the
# presentation is touched/re-saved with those values confirmed, without
# altering any of the effective transition semantics (same speed, same
# advance-on-click/time, same durations, same morph/fade effects).

$p = $ppt.ActivePresentation

$targetSlides = 10, 11, 12, 13

foreach ($idx in $targetSlides) {
    $s = $p.Slides.Item($idx)
    $t = $s.SlideShowTransition

    # Read back the current (unchanged) transition configuration. These
    # values are identical before and after the commit - only the XML
    # namespace plumbing around the compatibility wrapper changes.
    $speed          = $t.Speed
    $duration       = $t.Duration
    $advanceOnClick = $t.AdvanceOnClick
    $advanceOnTime  = $t.AdvanceOnTime
    $advanceTime    = $t.AdvanceTime

    Write-Output ("Slide " + $idx + ": speed=" + $speed + " duration=" + $duration + " advanceOnClick=" + $advanceOnClick + " advanceOnTime=" + $advanceOnTime + " advanceTime=" + $advanceTime)
}
